$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 9627879873
$ws.Range("C5").Value = "Saumya Tiwari"
$ws.Range("D5").Value = "Nai Basti"
$ws.Range("E5").Value = "House no. 8"
$ws.Range("F5").Value = "29.080523,80.110608"

$ws.Range("F5").Select()
